$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.306.67"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "2.082.36"
$ws.Range("E3").Value = "  +3.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.10"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5193"
$ws.Range("E7").Value = "  +1.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4311"
$ws.Range("E8").Value = "  +3.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08824"
$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.10"
$ws.Range("E10").Value = "  +6.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.161"
$ws.Range("E11").Value = "  +2.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.65"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").Value = "2.078.51"
$ws.Range("E13").Value = "  +2.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.682"
$ws.Range("E14").Value = "  +1.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.693"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.16"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06615"
$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9983"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.326"
$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("D23").Value = "30.350.62"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.30"
$ws.Range("E24").Value = "  +3.45%  "

$ws.Range("E25").Value = "  +2.19%  "

$ws.Range("D26").Value = "2.320.20"
$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.34"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.597"
$ws.Range("E28").Value = "  +6.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.73"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.99"
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.188"
$ws.Range("E31").Value = "  +4.13%  "

$ws.Range("E32").Value = "  +1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.628"
$ws.Range("E33").Value = "  +19.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.235"
$ws.Range("E34").Value = "  +2.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.821"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02584"
$ws.Range("E36").Value = "  +2.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.796"
$ws.Range("E37").Value = "  +7.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.68"
$ws.Range("E38").Value = "  +2.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06659"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.438"
$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2254"
$ws.Range("E41").Value = "  +2.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6850"
$ws.Range("E42").Value = "  +2.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9982"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.03"
$ws.Range("E45").Value = "  +2.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6362"
$ws.Range("E46").Value = "  +2.70%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.607"
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.243"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.191"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.61"
$ws.Range("E51").Value = "  +0.55%  "
